$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header "Semester" -> "SEMESTER"
$ws.Range("E1").Value = "SEMESTER"

# Update the active selection to H5 (cosmetic, matches the saved view state)
$ws.Range("H5").Select()
